# Commit: "Solved 252 in java"
# Adds a new row to the LeetCode tracker table for
# "252 - Meeting Rooms", inserted right after problem 205
# (Isomorphic Strings) and before problem 2 (Add Two Numbers),
# i.e. at worksheet row 44 - shifting every following row down
# by one and growing Table1 from A2:X74 to A2:X75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 44; rows 44..74 shift down to 45..75,
# and the new row inherits formatting from the row above it.
$ws.Rows.Item(44).Insert()

# Grow Table1 (and its autofilter) to cover the newly inserted row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:X75"))

# Fill in the new row's data.
$ws.Range("A44").Value = "Array"
$ws.Range("B44").Value = 252
$ws.Range("C44").Value = "252 - Meeting Rooms"
$ws.Range("D44").Value = "Easy"
$ws.Range("E44").Value = "Sort, then iterate"
$ws.Range("F44").Value = "O(nlogn) time, O(n) memory"
$ws.Range("G44").Value = "O(nlogn) time"
$ws.Range("I44").Value = "O(n) memory"
$ws.Range("J44").Value = "no"
$ws.Range("K44").Value = "no"
$ws.Range("M44").Value = "35 minutes"
